# Auto-generated edit script applying the Zalera_Profits sheet value updates
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 32
$ws.Range("H32").Value = 13112.5
$ws.Range("I32").Value = 17200
$ws.Range("J32").Value = 850
$ws.Range("K32").Value = 17200
$ws.Range("L32").Value = 850
$ws.Range("M32").Value = -16874
$ws.Range("N32").Value = -1502
# row 34
$ws.Range("H34").Value = 8466.166999999999
$ws.Range("I34").Value = 8466.166999999999
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 8466.166999999999
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -8263.166999999999
$ws.Range("N34").ClearContents()
# row 36
$ws.Range("H36").Value = 8466.166999999999
$ws.Range("I36").Value = 8466.166999999999
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 8466.166999999999
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -7751.166999999999
$ws.Range("N36").ClearContents()
# row 74
$ws.Range("H74").Value = 8300.200000000001
$ws.Range("I74").Value = 5667.6665
$ws.Range("J74").Value = 9428.429
$ws.Range("K74").Value = 5667.6665
$ws.Range("L74").Value = 9428.429
$ws.Range("M74").Value = -4731.6665
$ws.Range("N74").Value = -11300.429
# row 77
$ws.Range("H77").Value = 8300.200000000001
$ws.Range("I77").Value = 5667.6665
$ws.Range("J77").Value = 9428.429
$ws.Range("K77").Value = 28338.3325
$ws.Range("L77").Value = 47142.145
$ws.Range("M77").Value = -23658.3325
$ws.Range("N77").Value = -56502.145
# row 87
$ws.Range("H87").Value = 62352.234
$ws.Range("J87").Value = 62352.234
$ws.Range("L87").Value = 62352.234
$ws.Range("N87").Value = -64848.234
# row 90
$ws.Range("H90").Value = 62352.234
$ws.Range("J90").Value = 62352.234
$ws.Range("L90").Value = 187056.702
$ws.Range("N90").Value = -199536.702
# row 127
$ws.Range("H127").Value = 999.3125
$ws.Range("I127").Value = 732.6
$ws.Range("J127").Value = 5000
$ws.Range("K127").Value = 2197.8
$ws.Range("L127").Value = 15000
$ws.Range("M127").Value = 2762.2
$ws.Range("N127").Value = -24920
# row 129
$ws.Range("H129").Value = 2181.7646
$ws.Range("I129").Value = 1591.1111
$ws.Range("J129").Value = 2846.25
$ws.Range("K129").Value = 4773.3333
$ws.Range("L129").Value = 8538.75
$ws.Range("M129").Value = 226.6666999999998
$ws.Range("N129").Value = -18538.75
# row 135
$ws.Range("H135").Value = 3425.8823
$ws.Range("I135").Value = 5204.4
$ws.Range("J135").Value = 2684.8333
$ws.Range("K135").Value = 46839.6
$ws.Range("L135").Value = 24163.4997
$ws.Range("M135").Value = -44304.6
$ws.Range("N135").Value = -29233.4997
# row 138
$ws.Range("H138").Value = 4584.0713
$ws.Range("I138").Value = 4168.143
$ws.Range("K138").Value = 12504.429
$ws.Range("M138").Value = -7364.429

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 23209.982
$ws.Range("I32").Value = 22948.885
$ws.Range("K32").Value = 22948.885
$ws.Range("M32").Value = -22661.885
# row 45
$ws.Range("H45").Value = 2074.6667
$ws.Range("I45").Value = 1194.75
$ws.Range("K45").Value = 1194.75
$ws.Range("M45").Value = -817.75
# row 97
$ws.Range("H97").Value = 1280986.6
$ws.Range("I97").Value = 1687091.1
$ws.Range("J97").Value = 4658.143
$ws.Range("K97").Value = 1687091.1
$ws.Range("L97").Value = 4658.143
$ws.Range("M97").Value = -1686595.1
$ws.Range("N97").Value = -5650.143
# row 110
$ws.Range("H110").Value = 9616892
$ws.Range("I110").Value = 15625604
$ws.Range("K110").Value = 15625604
$ws.Range("M110").Value = -15623559
# row 122
$ws.Range("H122").Value = 2827.3
$ws.Range("I122").Value = 1896.1538
$ws.Range("K122").Value = 5688.4614
$ws.Range("M122").Value = -3238.4614

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 20
$ws.Range("H20").Value = 3843.9614
$ws.Range("I20").Value = 3783.261
$ws.Range("J20").Value = 4309.3335
$ws.Range("K20").Value = 3783.261
$ws.Range("L20").Value = 4309.3335
$ws.Range("M20").Value = -3536.261
$ws.Range("N20").Value = -4803.3335

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 43484896
$ws.Range("I31").Value = 200002240
$ws.Range("K31").Value = 200002240
$ws.Range("M31").Value = -200001945
# row 34
$ws.Range("H34").Value = 43484896
$ws.Range("I34").Value = 200002240
$ws.Range("K34").Value = 200002240
$ws.Range("M34").Value = -200002038
# row 122
$ws.Range("H122").Value = 64031.188
$ws.Range("I122").Value = 72463.21000000001
$ws.Range("J122").Value = 5007
$ws.Range("K122").Value = 217389.63
$ws.Range("L122").Value = 15021
$ws.Range("M122").Value = -214939.63
$ws.Range("N122").Value = -19921
# row 132
$ws.Range("H132").Value = 5159
$ws.Range("I132").Value = 2698.75
$ws.Range("K132").Value = 8096.25
$ws.Range("M132").Value = -5566.25
# row 133
$ws.Range("H133").Value = 100999
$ws.Range("J133").Value = 100999
$ws.Range("L133").Value = 100999
$ws.Range("N133").Value = -106059

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 3
$ws.Range("H3").Value = 4920.5713
$ws.Range("I3").Value = 2836
$ws.Range("K3").Value = 8508
$ws.Range("M3").Value = -8396
# row 37
$ws.Range("H37").Value = 79769
$ws.Range("J37").Value = 79769
$ws.Range("L37").Value = 239307
$ws.Range("N37").Value = -239531
# row 122
$ws.Range("H122").Value = 10819.286
$ws.Range("J122").Value = 1052.5
$ws.Range("L122").Value = 9472.5
$ws.Range("N122").Value = -14372.5
# row 131
$ws.Range("H131").Value = 3843.38
$ws.Range("I131").Value = 621.4
$ws.Range("K131").Value = 1864.2
$ws.Range("M131").Value = 3175.8
# row 132
$ws.Range("H132").Value = 40321.54
$ws.Range("I132").Value = 78103.69500000001
$ws.Range("K132").Value = 702933.2550000001
$ws.Range("M132").Value = -700403.2550000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 11
$ws.Range("H11").Value = 1887692.1
$ws.Range("J11").Value = 313000.28
$ws.Range("L11").Value = 313000.28
$ws.Range("N11").Value = -313278.28
# row 12
$ws.Range("H12").Value = 6499.6665
$ws.Range("J12").Value = 8999.5
$ws.Range("L12").Value = 8999.5
$ws.Range("N12").Value = -9279.5
# row 21
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
# row 30
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
# row 102
$ws.Range("H102").Value = 3465.5557
$ws.Range("I102").Value = 3741.4285
$ws.Range("K102").Value = 3741.4285
$ws.Range("M102").Value = -2119.4285
# row 138
$ws.Range("H138").Value = 89993
$ws.Range("J138").Value = 89993
$ws.Range("L138").Value = 89993
$ws.Range("N138").Value = -100273

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 253499.75
$ws.Range("I7").Value = 502000
$ws.Range("J7").Value = 4999.5
$ws.Range("K7").Value = 502000
$ws.Range("L7").Value = 4999.5
$ws.Range("M7").Value = -501888
$ws.Range("N7").Value = -5223.5
# row 46
$ws.Range("H46").Value = 5096.811
$ws.Range("I46").Value = 1134.1666
$ws.Range("K46").Value = 1134.1666
$ws.Range("M46").Value = -946.1666
# row 68
$ws.Range("H68").Value = 5347.5
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
# row 71
$ws.Range("H71").Value = 5347.5
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
# row 126
$ws.Range("H126").Value = 253499.75
$ws.Range("I126").Value = 502000
$ws.Range("J126").Value = 4999.5
$ws.Range("K126").Value = 1506000
$ws.Range("L126").Value = 14998.5
$ws.Range("M126").Value = -1503530
$ws.Range("N126").Value = -19938.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 81
$ws.Range("H81").Value = 12759.125
$ws.Range("J81").Value = 16970.559
$ws.Range("L81").Value = 33941.118
$ws.Range("N81").Value = -36063.118
# row 84
$ws.Range("H84").Value = 12759.125
$ws.Range("J84").Value = 16970.559
$ws.Range("L84").Value = 169705.59
$ws.Range("N84").Value = -180313.59
# row 107
$ws.Range("H107").Value = 1426.2
$ws.Range("I107").Value = 761.1667
$ws.Range("K107").Value = 2283.5001
$ws.Range("M107").Value = -363.5001000000002
# row 122
$ws.Range("H122").Value = 2806.8667
$ws.Range("I122").Value = 2864.5
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 8593.5
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -6143.5
$ws.Range("N122").Value = -10900
# row 132
$ws.Range("H132").Value = 5916.744
$ws.Range("I132").Value = 5146.5137
$ws.Range("J132").Value = 10666.5
$ws.Range("K132").Value = 15439.5411
$ws.Range("L132").Value = 31999.5
$ws.Range("M132").Value = -12909.5411
$ws.Range("N132").Value = -37059.5
# row 136
$ws.Range("H136").Value = 3551.077
$ws.Range("I136").Value = 1393.8
$ws.Range("J136").Value = 6492.8184
$ws.Range("K136").Value = 4181.4
$ws.Range("L136").Value = 19478.4552
$ws.Range("M136").Value = -1631.4
$ws.Range("N136").Value = -24578.4552
